$d = $word.ActiveDocument

# The five list-item sentences that need a trailing period appended,
# matched by their exact (paragraph-mark-trimmed) paragraph text.
$targets = @(
  "Used to store primitives ",
  "Each thread has its own stack memory",
  "All the threads share a same heap memory",
  "When we create object, it is created in Eden",
  "Is not part of heap memory"
)

foreach ($t in $targets) {
    $found = $false
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs($i)
        $text = $para.Range.Text
        # Strip the trailing paragraph mark (and any cell mark) so we can
        # compare against the plain sentence text.
        $trimmed = $text.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $t) {
            $r = $para.Range
            # Exclude the paragraph mark from the range, then drop a new
            # run containing just "." right after the existing text.
            [void]$r.MoveEnd(1, -1)
            [void]$r.Collapse(0)
            [void]$r.InsertAfter(".")
            $found = $true
            break
        }
    }
}
